$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("688:688").Insert()
$ws.Range("A688").NumberFormat = "@"
$ws.Range("A688").Value = "2026/01/20"
$ws.Range("A688").ClearFormats()
$ws.Range("B688").Value = "火"
$ws.Range("C688").Value = 13
$ws.Range("D688").Value = 171
